$d = $word.ActiveDocument

# The Title, Author and Abstract paragraphs each had their text split
# across many single-word runs (one run per word/space). Collapse each
# of those paragraphs down to a single run carrying the full text, by
# running a Find/Replace scoped to that paragraph's own range (old text
# == new text, so only the run structure changes, not the wording).

# Title paragraph -> "Questions: Solving exponential equations"
$titlePara = $d.Paragraphs(1).Range
[void]$titlePara.Find.Execute("Questions: Solving exponential equations", $true, $false, $false, $false, $false, $true, 1, $false, "Questions: Solving exponential equations", 2)

# Author paragraph -> "Zoë Gemmell, Isabella Lewis, Akshat Srivastava"
$authorPara = $d.Paragraphs(2).Range
[void]$authorPara.Find.Execute("Zoë Gemmell, Isabella Lewis, Akshat Srivastava", $true, $false, $false, $false, $false, $true, 1, $false, "Zoë Gemmell, Isabella Lewis, Akshat Srivastava", 2)

# Abstract paragraph -> "A selection of questions for the study guide on solving equations involving indices."
$abstractPara = $d.Paragraphs(4).Range
[void]$abstractPara.Find.Execute("A selection of questions for the study guide on solving equations involving indices.", $true, $false, $false, $false, $false, $true, 1, $false, "A selection of questions for the study guide on solving equations involving indices.", 2)
